# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 08:33"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1237761
$ws.Range("C4").Value = 128
$ws.Range("D4").Value = 200669
$ws.Range("E4").Value = 964817
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 72275

# Row 60 - Kazajistan
$ws.Range("B60").Value = 4277
$ws.Range("C60").Value = 72
$ws.Range("E60").Value = 2969

# Row 112 - Georgia
$ws.Range("B112").Value = 610
$ws.Range("C112").Value = 6
$ws.Range("D112").Value = 269
$ws.Range("E112").Value = 332

# Row 122 - Taiwan
$ws.Range("B122").Value = 439
$ws.Range("C122").Value = 1
$ws.Range("D122").Value = 339
$ws.Range("E122").Value = 94
